# Comentarios de Rosa a la prelectura
#
# Adds the four label textboxes ("a)", "b)", "c)", "d)") that Rosa's review
# comments call out on the single appendix slide. Positions/sizes below are
# the exact EMU values from the target OOXML, expressed in points
# (1 pt = 12700 EMU) because Shapes.AddTextbox takes its geometry in points.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# a)
$shp1 = $s.Shapes.AddTextbox(1, 275.6699212598425, 7.733228346456693, 45.359370078740156, 29.081259842519685)
$shp1.Name = "1 CuadroTexto"
$shp1.Fill.Visible = 0
$shp1.TextFrame.WordWrap = 1
$shp1.TextFrame.AutoSize = 1
$shp1.TextFrame.TextRange.Text = "a)"
$shp1.TextFrame.TextRange.LanguageID = "es-ES"

# b)
$shp2 = $s.Shapes.AddTextbox(1, 475.448031496063, 169.80685039370078, 129.10385826771653, 29.081259842519685)
$shp2.Name = "2 CuadroTexto"
$shp2.Fill.Visible = 0
$shp2.TextFrame.WordWrap = 1
$shp2.TextFrame.AutoSize = 1
$shp2.TextFrame.TextRange.Text = "b)"
$shp2.TextFrame.TextRange.LanguageID = "es-ES"

# c)
$shp3 = $s.Shapes.AddTextbox(1, 428.7577952755905, 473.88763779527557, 39.68944881889764, 29.081259842519685)
$shp3.Name = "3 CuadroTexto"
$shp3.Fill.Visible = 0
$shp3.TextFrame.WordWrap = 1
$shp3.TextFrame.AutoSize = 1
$shp3.TextFrame.TextRange.Text = "c)"
$shp3.TextFrame.TextRange.LanguageID = "es-ES"

# d) -- this one does not wrap (wrap="none" in the target XML)
$shp4 = $s.Shapes.AddTextbox(1, 139.81740157480314, 395.5903937007874, 29.68708661417323, 29.081259842519685)
$shp4.Name = "4 CuadroTexto"
$shp4.Fill.Visible = 0
$shp4.TextFrame.WordWrap = 0
$shp4.TextFrame.AutoSize = 1
$shp4.TextFrame.TextRange.Text = "d)"
$shp4.TextFrame.TextRange.LanguageID = "es-ES"
